# Apply edits described by the diff:
# 1. Cell D2 on sheet "general_report" changes from the numeric value 4711
#    to the text "3 - leicht umgehbar" (matching D3/D4 in the same column).
# 2. The active selection on the sheet moves from J2 to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general_report")

# Update D2 to hold the same text value as D3/D4 ("3 - leicht umgehbar")
$ws.Range("D2").Value = "3 - leicht umgehbar"

# Move the active selection to A2
$ws.Range("A2").Select()
